# Weekly update: a new price record (week of 2023-06-20) is prepended to the
# "Feria Lagunitas de Puerto Montt - Membrillo" table, just above the record
# that used to be the first row (old row 106), pushing every following row
# down by one. This grows the used range from A1:T166 to A1:T167.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at row 106 - everything that was on rows 106..166
# shifts down to 107..167.
$ws.Rows.Item(106).Insert()

# The row immediately below (107) now holds the data that used to live on
# row 106; reuse its unchanged descriptive columns (market/region/product
# metadata, quality grade, unit label, origin, kg-per-unit) for the new row.
$ws.Cells.Item(106, 1).Value2 = $ws.Cells.Item(107, 1).Value2   # Mercado ID
$ws.Cells.Item(106, 2).Value2 = $ws.Cells.Item(107, 2).Value2   # Mercado
$ws.Cells.Item(106, 3).Value2 = $ws.Cells.Item(107, 3).Value2   # Región
$ws.Cells.Item(106, 4).Value2 = 45097                            # Fecha
$ws.Cells.Item(106, 5).Value2 = $ws.Cells.Item(107, 5).Value2   # Codreg
$ws.Cells.Item(106, 6).Value2 = $ws.Cells.Item(107, 6).Value2   # Tipo
$ws.Cells.Item(106, 7).Value2 = $ws.Cells.Item(107, 7).Value2   # Producto ID
$ws.Cells.Item(106, 8).Value2 = $ws.Cells.Item(107, 8).Value2   # Producto
$ws.Cells.Item(106, 9).Value2 = $ws.Cells.Item(107, 9).Value2   # Categoría ID
$ws.Cells.Item(106, 10).Value2 = $ws.Cells.Item(107, 10).Value2 # Categoría
$ws.Cells.Item(106, 11).Value2 = $ws.Cells.Item(107, 11).Value2 # Variedad
$ws.Cells.Item(106, 12).Value2 = $ws.Cells.Item(107, 12).Value2 # Calidad
$ws.Cells.Item(106, 13).Value2 = 200                              # Volumen
$ws.Cells.Item(106, 14).Value2 = 12000                            # Precio mínimo
$ws.Cells.Item(106, 15).Value2 = 13000                            # Precio máximo
$ws.Cells.Item(106, 16).Value2 = 12500                            # Precio promedio ponderado
$ws.Cells.Item(106, 17).Value2 = $ws.Cells.Item(107, 17).Value2 # Unidad de comercialización
$ws.Cells.Item(106, 18).Value2 = $ws.Cells.Item(107, 18).Value2 # Origen
$ws.Cells.Item(106, 19).Value2 = 694                              # Precio $/Kg
$ws.Cells.Item(106, 20).Value2 = $ws.Cells.Item(107, 20).Value2 # Kg / unidad
